# "sentencia 2 git status. segunda rama creada"
# Appends two new slides (ids 257/258) after the existing title slide (id 256),
# each using the "Title and Content" layout (ppLayoutObject = 2), which maps
# to slideLayout2.xml: a `title` placeholder plus a content placeholder at
# idx="1" - matching the shape layout used by the new slides in the diff.

$p = $ppt.ActivePresentation

# --- Slide 2: "git init" ---------------------------------------------------
$slide2 = $p.Slides.Add(2, 2)
$slide2.Shapes.Item(1).TextFrame.TextRange.Text = "git init"
$slide2.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2  # ppAlignCenter
$slide2.Shapes.Item(2).TextFrame.TextRange.Text = "Crea un repositorio Git vacío o reinicializa uno existente"

# --- Slide 3: "Git status" --------------------------------------------------
$slide3 = $p.Slides.Add(3, 2)
$slide3.Shapes.Item(1).TextFrame.TextRange.Text = "Git status"
$slide3.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2  # ppAlignCenter
$slide3.Shapes.Item(2).TextFrame.TextRange.Text = "te mostrará los diferentes estados de los archivos en tu directorio de trabajo y área de ensayo. Qué archivos están modificados y sin seguimiento y cuáles con seguimiento pero no confirmados aún. "
